# Updated queries for C3DC first half testcases.
#
# The StatQuery / TabQuery cells on the ParticipantsTab/DiagnosisTab/
# TreatmentTab/TreatmentRespTab/SurvivalTab/StudiesTab rows all join on the
# old generic "id" columns (std.id / prt.id) -- these are updated to use the
# fully-qualified id columns (std.study_id / prt.participant_id) that match
# the renamed columns in the source dataframes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Update-C3DCQuery {
    param([string]$CellRef)

    $text = $ws.Range($CellRef).Value()
    if ($null -eq $text) { return }

    $text = $text -replace 'df_participant prt ON std\.id = prt\."study\.id"', 'df_participant prt ON std.study_id = prt."study.study_id"'
    $text = $text -replace 'df_diagnoses dgn ON prt\.id = dgn\."participant\.id"', 'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'
    $text = $text -replace 'df_treatments trt ON prt\.id = trt\."participant\.id"', 'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'
    $text = $text -replace 'df_treatment_resp trr ON prt\.id = trr\."participant\.id"', 'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'
    $text = $text -replace 'df_survival srv ON prt\.id = srv\."participant\.id"', 'df_survival srv ON prt.participant_id = srv."participant.participant_id"'
    $text = $text -replace 'df_reference_files rfs ON std\.id = rfs\."study\.id"', 'df_reference_files rfs ON std.study_id = rfs."study.study_id"'

    $ws.Range($CellRef).Value = $text
}

# StudiesTab row: StatQuery (C2) + TabQuery (B2)
Update-C3DCQuery "C2"
Update-C3DCQuery "B2"

# ParticipantsTab row: TabQuery (B3)
Update-C3DCQuery "B3"

# DiagnosisTab row: TabQuery (B4)
Update-C3DCQuery "B4"

# TreatmentTab row: TabQuery (B5)
Update-C3DCQuery "B5"

# TreatmentRespTab row: TabQuery (B6)
Update-C3DCQuery "B6"

# SurvivalTab row: TabQuery (B7)
Update-C3DCQuery "B7"

# Best-effort: restore the scrolled-to-top-left view (A3) shown in the
# updated workbook; harmless if the host doesn't persist it.
try {
    $excel.ActiveWindow.ScrollRow = 3
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
